$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$hf = $m.HeadersFooters
$dt = $hf.DateAndTime
$dt.Value = "30/06/2014"
Write-Output ("Value now: " + $dt.Value)
